$wb = $excel.ActiveWorkbook

# --- color sheet: zoom + selection change (no data changes) ---
$color = $wb.Worksheets.Item("color")
$color.Activate()
$excel.ActiveWindow.Zoom = 153
$color.Range("C1:C6").Select()

# --- person sheet: drop tabSelected, change selection, dedupe duplicate style (rows 8-9) ---
$person = $wb.Worksheets.Item("person")
$person.Activate()
$person.Range("A8:D9").HorizontalAlignment = -4108
$person.Range("C1:C8").Select()

# --- nation sheet: selection change only ---
$nation = $wb.Worksheets.Item("nation")
$nation.Activate()
$nation.Range("C1:C7").Select()

# --- pair sheet: "neuf" -> "nouveau", selection change ---
$pair = $wb.Worksheets.Item("pair")
$pair.Activate()
$pair.Range("E7").Value = "nouveau"
$pair.Range("G1:G8").Select()

# --- mood sheet: becomes active tab; move B2's text into C2 and D2; dedupe duplicate style ---
$mood = $wb.Worksheets.Item("mood")
$mood.Activate()
$moodB2 = $mood.Range("B2").Value2
$mood.Range("C2").Value = $moodB2
$mood.Range("D2").Value = $moodB2
$mood.Range("A2").HorizontalAlignment = -4108
$mood.Range("C2:D2").HorizontalAlignment = -4108
$mood.Range("B2").Clear()
$mood.Range("C4").HorizontalAlignment = -4108
$mood.Range("A5:D6").HorizontalAlignment = -4108
$mood.Range("A8:D10").HorizontalAlignment = -4108
$mood.Range("D3:D4").Select()

Write-Host "done"
